# Edit: rename/delete Cliona bowtie sheets, add Pione clean/un columns to "pione bowtie"
$wb = $excel.ActiveWorkbook

# --- 1. Rename / delete worksheets -----------------------------------------
$wb.Worksheets.Item("bowtie").Name  = "cliona bowtie"
$wb.Worksheets.Item("bowtie2").Name = "cliona bowtie2"
$wb.Worksheets.Item("bowtie3").Delete()
$wb.Worksheets.Item("bowtie4").Name = "cliona bowtie3"
$wb.Worksheets.Item("bowtie5").Name = "cliona bowtie4"

# --- 2. "pione bowtie" sheet: add F (clean) / H (un) columns + new J formula
$ws = $wb.Worksheets.Item("pione bowtie")

for ($row = 1; $row -le 22; $row++) {
    $a = $ws.Cells.Item($row, 1).Value2
    $base = $a -replace '\.sam$', ''
    $ws.Cells.Item($row, 6).Value = "$base.clean"
    $ws.Cells.Item($row, 8).Value = "$base.un"
}

$ws.Range("J1").Formula = '=("bowtie2 --local -x /mnt/beegfs/home/mstudiva/db/Plampa -1 "&C1&" -2 "&C23&" -S "&A1&" --no-hd --no-sq --no-unal --al-conc ./"&F1&" --un-conc junk/"&H1)'
$ws.Range("J2:J22").Formula = '=("bowtie2 --local -x /mnt/beegfs/home/mstudiva/db/Plampa -1 "&C2&" -2 "&C24&" -S "&A2&" --no-hd --no-sq --no-unal --al-conc ./"&F2&" --un-conc junk/"&H2)'

# Column widths for the new F/H columns (auto-fit to content)
$ws.Columns.Item(6).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null

Write-Host "done"
